$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4, 5, 6 have their species/location data cyclically rotated:
#   new row4 = old row6
#   new row5 = old row4
#   new row6 = old row5
# Capture old values first, then write them to their new positions.

$cols = @("A","B","E","F","G","H","P","Q","R","S","AI")

$old4 = @{}
$old5 = @{}
$old6 = @{}
foreach ($col in $cols) {
    $old4[$col] = $ws.Range("$col`4").Value2
    $old5[$col] = $ws.Range("$col`5").Value2
    $old6[$col] = $ws.Range("$col`6").Value2
}

foreach ($col in $cols) {
    $ws.Range("$col`4").Value2 = $old6[$col]
    $ws.Range("$col`5").Value2 = $old4[$col]
    $ws.Range("$col`6").Value2 = $old5[$col]
}
